# Regenerate orders with updated distance/size codes.
# Distance codes: D80->D86, D64->D69, D51->D55
# Size code:      S30->S31
# These substrings appear throughout Condition, Filename_Left,
# Filename_Right, Distance and Size columns (and their shared-string
# derived combinations), so a straightforward find/replace across the
# whole used range reproduces the regenerated order file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("S30", "S31")
